$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")
$ws.Columns.Item(3).Delete()
